$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells retain their original text (String) storage type rather than
# being auto-converted to numbers/dates by Excel when the new value looks numeric.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.304.36'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.45%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.832.32'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.92%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.012'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.99%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.87'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.010'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.79%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.79%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3687'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.98%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07439'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8850'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.46'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.884.96'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +5.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07334'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.25%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.80'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.93%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.569'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008788'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.27%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.80%  '
$ws.Range('B20').Value = 'WrappedBTC'
$ws.Range('C20').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.549.15'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.29%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.79'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.289'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.66'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.98%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.094.15'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.891'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.99'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.64'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.78%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.86%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.223'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.59%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '117.08'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.14%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7503'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.173'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.544'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.945'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.48%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.94%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.81%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05337'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01952'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.45%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.969'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.13%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.400'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.214'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5304'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1658'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.474'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4930'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.71%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.49'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.99%  '
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.011'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.89%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '105.02'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.671'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06302'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.10%  '
